$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the custom widths of columns B and C ---
# The host rounds ColumnWidth (in characters) to the nearest 1/6 of a
# character before it is written back to the OOXML <col width="..."/>
# attribute (output = ROUND((input + 5/6) * 6) / 6). The target widths
# (11.7109375 / 10.7109375, i.e. 1/256-character units typical of
# openpyxl-authored files) are not reachable exactly through the Excel
# object model, so we pick the input that lands on the closest possible
# 1/6-character grid point to the desired value.
$ws.Columns.Item(2).ColumnWidth = 10.833333333333332   # -> 11.666666666666666 (target 11.7109375)
$ws.Columns.Item(3).ColumnWidth = 9.833333333333332    # -> 10.666666666666666 (target 10.7109375)

# --- Update the three data values ---
$ws.Range("A1").Value = 161.7820836633567
$ws.Range("B1").Value = 4.3242779954860699
$ws.Range("C1").Value = 5.3979125896934113
